# Reorder glossary rows 2-17 per updated terms_definition_2 export.
# Row content is unchanged; only the row order is updated (see mapping below).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 = original row 16 (Uncertainty)
$ws.Range("A2").Value = 'Uncertainty'
$ws.Range("B2").Value = ''
$ws.Range("C2").Value = ''
$ws.Range("D2").Value = 'Non-negative parameter, associated with Data, which characterizes the dispersion of the values of a [Trait ]that could reasonably be attributed to a Phenomenon [by means of sensing or modelling].'
$ws.Range("E2").Value = '- In case of quantitative(continuous) Data the uncertainty may be, for example, a standard deviation (or a given multiple of it), or the half-width of an interval having a stated level of confidence. (see e.g. standard and Expanded uncertainty)
- For qualitative (categorical?) Data uncertainty may be, for example, expressed by commission and omission (‘confusion matrix’) or overall errors.'
$ws.Range("F2").Value = ''
$ws.Range("G2").Value = '- modified from GUM, VIM4 :3.1, FIDUCEO, Notes added'

# New row 3 = original row 15 (Verification)
$ws.Range("A3").Value = 'Verification'
$ws.Range("B3").Value = ''
$ws.Range("C3").Value = ''
$ws.Range("D3").Value = 'The evaluation of whether or not a product, service, or system complies with a regulation requirement, specification, or imposed condition. It is often an internal Process.'
$ws.Range("E3").Value = ''
$ws.Range("F3").Value = ''
$ws.Range("G3").Value = '- EU-US Land Imaging EO Collaboration'

# New row 4 = original row 3 (Ancillary Data)
$ws.Range("A4").Value = 'Ancillary Data'
$ws.Range("B4").Value = ''
$ws.Range("C4").Value = ''
$ws.Range("D4").Value = 'Data other than instrument measurements, originating in the instrument itself or from the satellite, required to perform processing of the Data. They include orbit Data, attitude Data, time Information, and spacecraft engineering Data, Calibration Data, Data quality Information, and Data from other instruments or earth system models.'
$ws.Range("E4").Value = ''
$ws.Range("F4").Value = ''
$ws.Range("G4").Value = '- CEOS-ARD PFS template 20220302'

# New row 5 = original row 2 (Geolocating)
$ws.Range("A5").Value = 'Geolocating'
$ws.Range("B5").Value = ''
$ws.Range("C5").Value = ''
$ws.Range("D5").Value = 'Determination of the geographic location of a >=2D feature(?).'
$ws.Range("E5").Value = ''
$ws.Range("F5").Value = ''
$ws.Range("G5").Value = '- ISO 19130-1:2018, 3.36 (‘geopositioning’), modified'

# New row 6 = original row 10 (Entity)
$ws.Range("A6").Value = 'Entity'
$ws.Range("B6").Value = '- base'
$ws.Range("C6").Value = ''
$ws.Range("D6").Value = 'A government or business organization that is formed to conduct business or represent the government of the day.'
$ws.Range("E6").Value = ''
$ws.Range("F6").Value = 'CEOS Entities include Working Groups, Virtual Constellations, etc.'
$ws.Range("G6").Value = '- WGISS Shared Collection Lifecycle Management Principles for Earth Observation Data)'

# New row 7 = original row 6 (Georectifying)
$ws.Range("A7").Value = 'Georectifying'
$ws.Range("B7").Value = ''
$ws.Range("C7").Value = 'Orthorectifying'
$ws.Range("D7").Value = 'The correction of sample locations to achieve some sort of geometric regularity, e.g., a regular 2D geographic grid.'
$ws.Range("E7").Value = ''
$ws.Range("F7").Value = ''
$ws.Range("G7").Value = '- KCEO'

# New row 8 = original row 9 (Characteristic)
$ws.Range("A8").Value = 'Characteristic'
$ws.Range("B8").Value = '- base'
$ws.Range("C8").Value = ''
$ws.Range("D8").Value = 'Abstraction of a Property of an Object or of a set of objects.'
$ws.Range("E8").Value = '- Characteristics are used for describing Concepts.'
$ws.Range("F8").Value = ''
$ws.Range("G8").Value = '- ISO 1087-1:2000, 3.2.4; ISO 19146:2010(E); https://www.iso.org/standard/20057.html'

# New row 9 = original row 13 (Reference)
$ws.Range("A9").Value = 'Reference'
$ws.Range("B9").Value = ''
$ws.Range("C9").Value = ''
$ws.Range("D9").Value = 'A sort of Data acquired with an Uncertainty significantly lower (quantify?) than that of the Data it is being compared with.'
$ws.Range("E9").Value = ''
$ws.Range("F9").Value = ''
$ws.Range("G9").Value = '- VIM?, modified'

# New row 10 = original row 4 (Instrument Data)
$ws.Range("A10").Value = 'Instrument Data'
$ws.Range("B10").Value = ''
$ws.Range("C10").Value = ''
$ws.Range("D10").Value = 'Data created by an instrument including scientific measurements and any engineering or ancillary data which may be included in the data packets.'
$ws.Range("E10").Value = ''
$ws.Range("F10").Value = ''
$ws.Range("G10").Value = '- [EO Data Stewardship Glossary](https://ceos.org/document_management/Working_Groups/WGISS/Interest_Groups/Data_Stewardship/White_Papers/EO-DataStewardshipGlossary.pdf)'

# New row 11 = original row 7 (Test term)
$ws.Range("A11").Value = 'Test term'
$ws.Range("B11").Value = ''
$ws.Range("C11").Value = ''
$ws.Range("D11").Value = 'Second_definition_goes_here.'
$ws.Range("E11").Value = '- here should be bullets
- like this'
$ws.Range("F11").Value = '- this is also bullets
- like this'
$ws.Range("G11").Value = '- KCEO (no link included, so no brackets)
- [Website](https://en.wikipedia.org/wiki/Thai_script) ( if you have web references, just add the term goes into square [] brackets and the url into () normal brackets
```

---

References: 

1. Strobl, P. A., Woolliams, E. R., & Molch, K. (2024). Lost in Translation: The Need for Common Vocabularies and an Interoperable Thesaurus in Earth Observation Sciences. Surveys in Geophysics, 1-29.'

# New row 12 = original row 5 (Quantity)
$ws.Range("A12").Value = 'Quantity'
$ws.Range("B12").Value = '- base'
$ws.Range("C12").Value = ''
$ws.Range("D12").Value = 'Property whose instances can be compared by ratio or only by order.'
$ws.Range("E12").Value = ''
$ws.Range("F12").Value = ''
$ws.Range("G12").Value = '- gEOGlos(VIM4 Notes omitted)'

# New row 13 = original row 11 (Baseline)
$ws.Range("A13").Value = 'Baseline'
$ws.Range("B13").Value = ''
$ws.Range("C13").Value = ''
$ws.Range("D13").Value = 'Source data that has been processed to a common set of requirements and organised into a form that allows immediate analysis and interoperability through time and with other collections.'
$ws.Range("E13").Value = ''
$ws.Range("F13").Value = ''
$ws.Range("G13").Value = '- WGISS Shared Collection Lifecycle Management Principles for Earth Observation Data)'

# New row 14 = original row 8 (Data)
$ws.Range("A14").Value = 'Data'
$ws.Range("B14").Value = '- core'
$ws.Range("C14").Value = ''
$ws.Range("D14").Value = 'Scientific or technical measurements, values calculated therefrom, observations, or facts that can be represented by numbers, tables, graphs, models, text, or symbols which are used as a basis for reasoning and further calculation.'
$ws.Range("E14").Value = ''
$ws.Range("F14").Value = ''
$ws.Range("G14").Value = '- WGISS Shared Collection Lifecycle Management Principles for Earth Observation Data)'

# New row 15 = original row 14 (Validation)
$ws.Range("A15").Value = 'Validation'
$ws.Range("B15").Value = ''
$ws.Range("C15").Value = ''
$ws.Range("D15").Value = 'Validation aims to verify that the specified requirements are achieved or compliant. This involves comparing  mission products with representative Reference Data, considering various Observation conditions, ensuring the quality and Traceability of the Reference Data used.'
$ws.Range("E15").Value = '- In this part of ISO 19159, the term validation is used in a limited sense and only relates to the validation of Calibration Data in order to control their change over time.'
$ws.Range("F15").Value = ''
$ws.Range("G15").Value = '- BIPM; QA4EO; ESA ?, modified'

# New row 16 = original row 17 (User)
$ws.Range("A16").Value = 'User'
$ws.Range("B16").Value = ''
$ws.Range("C16").Value = ''
$ws.Range("D16").Value = 'External person, institution or system that consumes provided services.'
$ws.Range("E16").Value = 'Includes Data Access or Science and Service Exploitation Platforms provided by a payload data ground segment.'
$ws.Range("F16").Value = ''
$ws.Range("G16").Value = '- EO Data Stewardship Glossary)'

# New row 17 = original row 12 (Auxiliary Data)
$ws.Range("A17").Value = 'Auxiliary Data'
$ws.Range("B17").Value = ''
$ws.Range("C17").Value = ''
$ws.Range("D17").Value = 'Data required to perform processing of Sensor Data which is not obtained from the Sensor itself. Include: (a) Data provided by the spacecraft (e.g. orbit Position and velocity, attitude, instrument house-keeping Data, on-board time), (b) Data not available from on-board sources.'
$ws.Range("E17").Value = 'For EnMAP, this includes (a) Orbit files, attitude files, Calibration Data, instrument house-keeping Data, (b) atmospheric parameters, Reference images.'
$ws.Range("F17").Value = ''
$ws.Range("G17").Value = '- ENMAP Glossary of Terms, https://www.enmap.org/Data/doc/EnMAP_Terms.pdf, 20210624
- EO Data Stewardship Glossary)'

